$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2418750077486038
$ws.Range("B1").Value = 0.1969166398048401
$ws.Range("C1").Value = 0.1748537868261337
$ws.Range("D1").Value = 0.1919311583042145
$ws.Range("E1").Value = 0.238029956817627
